$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.57830923414868
$ws.Range("C2").Value = 0.0344533211606062
$ws.Range("D2").Value = 45.8100752258772

$ws.Range("B3").Value = 0.0450645780754445
$ws.Range("C3").Value = 0.0436511093027488
$ws.Range("D3").Value = 1.03238105045378
$ws.Range("E3").Value = 0.301908608213814

$ws.Range("B4").Value = 0.0704838805608534
$ws.Range("C4").Value = 0.071404732994947
$ws.Range("D4").Value = 0.987103761957085
$ws.Range("E4").Value = 0.323606086640475
